$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.375.03'
$ws.Range('D3').Value = '1.882.39'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7118'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08013'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3169'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08344'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '1.901.10'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.268'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7186'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.369'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008652'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.03%  '
$ws.Range('D18').Value = '29.386.61'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '2.151.16'
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.834'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.104'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.511'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.444'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.353'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.206'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05411'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.948'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.188'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.688'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('D39').Value = '1.274.84'
$ws.Range('E39').Value = '  +3.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.751'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.523'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9182'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '113.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '74.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.79%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('E46').Value = '  +5.17%  '
$ws.Range('D47').Value = '2.045.12'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.817'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5224'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.568'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4383'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.24%  '
